# Refresh the "cryptos" price/volume snapshot to the values captured in the
# latest GitHub Actions run (coin name/link swap for rows 45-46 included).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.328.38"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "2.070.01"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.83"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.80"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0764"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "2.374.92"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.67"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.77"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.776"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "2.071.99"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "37.254.54"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.33"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.45"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "0.0₃0813"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.03"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.76"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.78"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.16"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("E30").Value = "  -4.26%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0616"
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("E39").Value = "  -4.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.42"
$ws.Range("E41").Value = "  +4.06%  "
$ws.Range("D42").Value = "1.480.90"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.29"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0938"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0212"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.07"
$ws.Range("E48").Value = "  -8.56%  "
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").Value = "2.264.04"
